$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A106").Value = 45953
$ws.Range("B106").Value = "四方坪站充电量(kw)"
$ws.Range("C106").Value = 667.08100000000013
$ws.Range("D106").Value = 1653.7469999999992
$ws.Range("E106").Value = 499.49999999999994
$ws.Range("F106").Value = 391.28300000000002
$ws.Range("G106").Value = 339.226
$ws.Range("H106").Value = 643.15899999999988
$ws.Range("I106").Value = 538.99900000000002
$ws.Range("J106").Value = 234.21199999999999
$ws.Range("K106").Value = 191.19099999999997
$ws.Range("L106").Value = 249.21700000000004
$ws.Range("M106").Value = 166.32000000000002
$ws.Range("N106").Value = 208.91399999999996
$ws.Range("O106").Value = 575.14199999999994
$ws.Range("P106").Value = 1932.2579999999998
$ws.Range("Q106").Value = 631.38200000000006
$ws.Range("R106").Value = 399.44099999999997
$ws.Range("S106").Value = 250.32699999999997
$ws.Range("T106").Value = 61.343000000000004
$ws.Range("U106").Value = 115.816
$ws.Range("V106").Value = 117.5
$ws.Range("W106").Value = 45.81
$ws.Range("X106").Value = 94.54
$ws.Range("Y106").Value = 0
$ws.Range("Z106").Value = 15.94

$ws.Range("A107").Value = 45953
$ws.Range("B107").Value = "高岭站充电量(kw)"
$ws.Range("C107").Value = 401.74799999999999
$ws.Range("D107").Value = 343.37299999999999
$ws.Range("E107").Value = 90.745000000000005
$ws.Range("F107").Value = 42.969000000000001
$ws.Range("G107").Value = 62.645000000000003
$ws.Range("H107").Value = 231.13
$ws.Range("I107").Value = 249.79
$ws.Range("J107").Value = 229.44800000000001
$ws.Range("K107").Value = 285.13299999999998
$ws.Range("L107").Value = 271.70599999999996
$ws.Range("M107").Value = 330.9129999999999
$ws.Range("N107").Value = 326.52199999999999
$ws.Range("O107").Value = 561.93600000000004
$ws.Range("P107").Value = 585.37100000000009
$ws.Range("Q107").Value = 282.61699999999996
$ws.Range("R107").Value = 271.65000000000003
$ws.Range("S107").Value = 390.19900000000007
$ws.Range("T107").Value = 193.51999999999998
$ws.Range("U107").Value = 1.7509999999999999
$ws.Range("V107").Value = 36.227000000000004
$ws.Range("W107").Value = 79.051000000000002
$ws.Range("X107").Value = 62.741999999999997
$ws.Range("Y107").Value = 20.72
$ws.Range("Z107").Value = 10.468999999999999

$ws.Range("F109").Select()
